$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New section header row (row 18): "1400 / 04"
$ws.Range("A18:C18").Merge()
$ws.Cells.Item(18, 1).Value = "1400 / 04"
$ws.Cells.Item(18, 1).Style = $ws.Cells.Item(14, 1).Style
$ws.Cells.Item(18, 2).Style = $ws.Cells.Item(14, 2).Style
$ws.Cells.Item(18, 3).Style = $ws.Cells.Item(14, 3).Style

# New data row (row 19): Report 10_4 / 2021 June 22 / 1400/04/01
$ws.Cells.Item(19, 1).Value = "Report 10_4"
$ws.Cells.Item(19, 2).Value = "2021 June 22"
$ws.Cells.Item(19, 3).Value = "1400/04/01"
$ws.Cells.Item(19, 1).Style = $ws.Cells.Item(13, 1).Style
$ws.Cells.Item(19, 2).Style = $ws.Cells.Item(13, 2).Style
$ws.Cells.Item(19, 3).Style = $ws.Cells.Item(13, 3).Style

# New blank row (row 20) with style 5 (like old row 19)
$ws.Cells.Item(20, 1).Style = $ws.Cells.Item(12, 1).Style
$ws.Cells.Item(20, 2).Style = $ws.Cells.Item(12, 2).Style
$ws.Cells.Item(20, 3).Style = $ws.Cells.Item(12, 3).Style

# New blank row (row 21) with style 1 (like old row 18/20)
$ws.Cells.Item(21, 1).Style = $ws.Cells.Item(11, 1).Style
$ws.Cells.Item(21, 2).Style = $ws.Cells.Item(11, 2).Style
$ws.Cells.Item(21, 3).Style = $ws.Cells.Item(11, 3).Style

$ws.Range("B17").Select()
